# Week 16 log + season sim from Week 17
# - Rushing sheet: insert a new row for G.Olszewski, rename J.Johnson -> Jak.Johnson,
#   and update weekly rushing totals.
# - Receiving sheet: rename J.Johnson -> Jak.Johnson and update weekly receiving totals.

$wb = $excel.ActiveWorkbook
$rushing = $wb.Worksheets.Item("Rushing")
$receiving = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------------
# Rushing sheet updates
# ---------------------------------------------------------------------------

# M.Jones (row 2)
$rushing.Range("C2").Value = 6
$rushing.Range("D2").Value = 7
$rushing.Range("E2").Value = 16

# D.Harris (row 3)
$rushing.Range("C3").Value = 113
$rushing.Range("D3").Value = 55
$rushing.Range("E3").Value = 15
$rushing.Range("F3").Value = 39

# B.Bolden (row 7)
$rushing.Range("D7").Value = 9
$rushing.Range("E7").Value = 18

# Insert a new row 11 for G.Olszewski (copy formatting from the row above it,
# then overwrite the copied values), pushing J.Johnson/J.Smith down a row.
$rushing.Range("A10:F10").Copy()
$rushing.Range("A11:F11").Insert()
$rushing.Range("A11").Borders.LineStyle = 1

$rushing.Range("A11").Value = 9
$rushing.Range("B11").Value = "G.Olszewski"
$rushing.Range("C11").Value = 1
$rushing.Range("D11").Value = 0
$rushing.Range("E11").Value = 0
$rushing.Range("F11").Value = 0

# J.Johnson (now row 12 after the insert) is renamed to Jak.Johnson; stats unchanged.
# The sequential index in column A shifts up by one for the rows pushed down.
$rushing.Range("A12").Value = 10
$rushing.Range("B12").Value = "Jak.Johnson"
$rushing.Range("A13").Value = 11

# ---------------------------------------------------------------------------
# Receiving sheet updates
# ---------------------------------------------------------------------------

# D.Harris (row 2)
$receiving.Range("C2").Value = 16

# J.Taylor (row 5)
$receiving.Range("C5").Value = 38
$receiving.Range("D5").Value = 31
$receiving.Range("G5").Value = 8
$receiving.Range("H5").Value = 8

# J.Meyers (row 7)
$receiving.Range("C7").Value = 86
$receiving.Range("D7").Value = 62
$receiving.Range("E7").Value = 25
$receiving.Range("F7").Value = 10
$receiving.Range("G7").Value = 12
$receiving.Range("H7").Value = 7

# K.Bourne (row 8)
$receiving.Range("C8").Value = 46
$receiving.Range("D8").Value = 38
$receiving.Range("E8").Value = 13
$receiving.Range("F8").Value = 9

# N.Harry (row 10)
$receiving.Range("C10").Value = 14
$receiving.Range("D10").Value = 9
$receiving.Range("E10").Value = 7

# J.Johnson (row 11) renamed to Jak.Johnson
$receiving.Range("B11").Value = "Jak.Johnson"
$receiving.Range("C11").Value = 4
$receiving.Range("D11").Value = 3

# H.Henry (row 13)
$receiving.Range("C13").Value = 53
$receiving.Range("D13").Value = 37
$receiving.Range("E13").Value = 12
